$wb = $excel.ActiveWorkbook

# --- Submit orders ---
$ws = $wb.Worksheets.Item("Submit orders")
$ws.Cells.Item(88, 1).Value = "10.21.2022 10:20 (Kyiv+Israel) 07:20 (UTC) 16:20 (Japan) 12:50 (India)"
$ws.Cells.Item(88, 2).Value = 1.583
$ws.Cells.Item(88, 3).Value = -0.856
$ws.Cells.Item(88, 4).Value = "***"
$ws.Cells.Item(88, 5).Value = "***"
$ws.Range($ws.Cells.Item(88, 1), $ws.Cells.Item(88, 5)).Style = "Normal"
$ws.Cells.Item(89, 1).Value = "10.21.2022 15:04 (Kyiv+Israel) 12:04 (UTC) 21:04 (Japan) 17:34 (India)"
$ws.Cells.Item(89, 2).Value = 0.862
$ws.Cells.Item(89, 3).Value = -0.135
$ws.Cells.Item(89, 4).Value = "***"
$ws.Cells.Item(89, 5).Value = "***"
$ws.Range($ws.Cells.Item(89, 1), $ws.Cells.Item(89, 5)).Style = "Normal"
$ws.Cells.Item(90, 1).Value = "10.24.2022 13:03 (Kyiv+Israel) 10:03 (UTC) 19:03 (Japan) 15:33 (India)"
$ws.Cells.Item(90, 2).Value = 1.02
$ws.Cells.Item(90, 3).Value = -0.293
$ws.Cells.Item(90, 4).Value = "***"
$ws.Cells.Item(90, 5).Value = "***"
$ws.Range($ws.Cells.Item(90, 1), $ws.Cells.Item(90, 5)).Style = "Normal"

# --- Submit internet survey ---
$ws = $wb.Worksheets.Item("Submit internet survey")
$ws.Cells.Item(81, 1).Value = "10.21.2022 11:15 (Kyiv+Israel) 08:15 (UTC) 17:15 (Japan) 13:45 (India)"
$ws.Cells.Item(81, 2).Value = 0.735
$ws.Cells.Item(81, 3).Value = -0.104
$ws.Cells.Item(81, 4).Value = "***"
$ws.Cells.Item(81, 5).Value = "***"
$ws.Range($ws.Cells.Item(81, 1), $ws.Cells.Item(81, 5)).Style = "Normal"
$ws.Cells.Item(82, 1).Value = "10.21.2022 15:07 (Kyiv+Israel) 12:07 (UTC) 21:07 (Japan) 17:37 (India)"
$ws.Cells.Item(82, 2).Value = 0.915
$ws.Cells.Item(82, 3).Value = -0.284
$ws.Cells.Item(82, 4).Value = "***"
$ws.Cells.Item(82, 5).Value = "***"
$ws.Range($ws.Cells.Item(82, 1), $ws.Cells.Item(82, 5)).Style = "Normal"
$ws.Cells.Item(83, 1).Value = "10.24.2022 13:06 (Kyiv+Israel) 10:06 (UTC) 19:06 (Japan) 15:36 (India)"
$ws.Cells.Item(83, 2).Value = 1.501
$ws.Cells.Item(83, 3).Value = -0.8699999999999999
$ws.Cells.Item(83, 4).Value = "***"
$ws.Cells.Item(83, 5).Value = "***"
$ws.Range($ws.Cells.Item(83, 1), $ws.Cells.Item(83, 5)).Style = "Normal"

# --- Submit a phone survey ---
$ws = $wb.Worksheets.Item("Submit a phone survey")
$ws.Cells.Item(77, 1).Value = "10.21.2022 11:29 (Kyiv+Israel) 08:29 (UTC) 17:29 (Japan) 13:59 (India)"
$ws.Cells.Item(77, 2).Value = 1.741
$ws.Cells.Item(77, 3).Value = -0.637
$ws.Cells.Item(77, 4).Value = "***"
$ws.Cells.Item(77, 5).Value = "***"
$ws.Range($ws.Cells.Item(77, 1), $ws.Cells.Item(77, 5)).Style = "Normal"
$ws.Cells.Item(78, 1).Value = "10.21.2022 11:55 (Kyiv+Israel) 08:55 (UTC) 17:55 (Japan) 14:25 (India)"
$ws.Cells.Item(78, 2).Value = 1.66
$ws.Cells.Item(78, 3).Value = -0.5559999999999998
$ws.Cells.Item(78, 4).Value = "***"
$ws.Cells.Item(78, 5).Value = "***"
$ws.Range($ws.Cells.Item(78, 1), $ws.Cells.Item(78, 5)).Style = "Normal"
$ws.Cells.Item(79, 1).Value = "10.21.2022 14:29 (Kyiv+Israel) 11:29 (UTC) 20:29 (Japan) 16:59 (India)"
$ws.Cells.Item(79, 2).Value = 1.42
$ws.Cells.Item(79, 3).Value = -0.3159999999999998
$ws.Cells.Item(79, 4).Value = "***"
$ws.Cells.Item(79, 5).Value = "***"
$ws.Range($ws.Cells.Item(79, 1), $ws.Cells.Item(79, 5)).Style = "Normal"
$ws.Cells.Item(80, 1).Value = "10.21.2022 15:10 (Kyiv+Israel) 12:10 (UTC) 21:10 (Japan) 17:40 (India)"
$ws.Cells.Item(80, 2).Value = 1.637
$ws.Cells.Item(80, 3).Value = -0.5329999999999999
$ws.Cells.Item(80, 4).Value = "***"
$ws.Cells.Item(80, 5).Value = "***"
$ws.Range($ws.Cells.Item(80, 1), $ws.Cells.Item(80, 5)).Style = "Normal"
$ws.Cells.Item(81, 1).Value = "10.24.2022 13:48 (Kyiv+Israel) 10:48 (UTC) 19:48 (Japan) 16:18 (India)"
$ws.Cells.Item(81, 2).Value = 2.375
$ws.Cells.Item(81, 3).Value = -1.271
$ws.Cells.Item(81, 4).Value = "***"
$ws.Cells.Item(81, 5).Value = "***"
$ws.Range($ws.Cells.Item(81, 1), $ws.Cells.Item(81, 5)).Style = "Normal"

# --- Checkertificate ---
$ws = $wb.Worksheets.Item("Checkertificate")
$ws.Cells.Item(92, 1).Value = "10.21.2022 11:24 (Kyiv+Israel) 08:24 (UTC) 17:24 (Japan) 13:54 (India)"
$ws.Cells.Item(92, 2).Value = 0.878
$ws.Cells.Item(92, 3).Value = -0.213
$ws.Cells.Item(92, 4).Value = "***"
$ws.Cells.Item(92, 5).Value = "***"
$ws.Range($ws.Cells.Item(92, 1), $ws.Cells.Item(92, 5)).Style = "Normal"
$ws.Cells.Item(93, 1).Value = "10.21.2022 15:19 (Kyiv+Israel) 12:19 (UTC) 21:19 (Japan) 17:49 (India)"
$ws.Cells.Item(93, 2).Value = 1.202
$ws.Cells.Item(93, 3).Value = -0.5369999999999999
$ws.Cells.Item(93, 4).Value = "***"
$ws.Cells.Item(93, 5).Value = "***"
$ws.Range($ws.Cells.Item(93, 1), $ws.Cells.Item(93, 5)).Style = "Normal"
$ws.Cells.Item(94, 1).Value = "10.24.2022 13:53 (Kyiv+Israel) 10:53 (UTC) 19:53 (Japan) 16:23 (India)"
$ws.Cells.Item(94, 2).Value = 1.16
$ws.Cells.Item(94, 3).Value = -0.4949999999999999
$ws.Cells.Item(94, 4).Value = "***"
$ws.Cells.Item(94, 5).Value = "***"
$ws.Range($ws.Cells.Item(94, 1), $ws.Cells.Item(94, 5)).Style = "Normal"
